$d = $word.ActiveDocument
$nl = [char]11

# --- Paragraphs whose entire run-text is replaced wholesale (single logical run) ---
$d.Paragraphs.Item(6).Range.Text = "Uma breve história sobre o uso de microrganismos na indústria farmacêutica, conceitos gerais sobre medicamentos biológicos, enzimas em medicamentos, proteínas terapêuticas, biologia molecular e sintética, expressão e produção de proteínas de interesse na indústria farmacêutica."
$d.Paragraphs.Item(7).Range.Text = "A brief history of the use of microorganisms in the pharmaceutical industry, general concepts about biological medicines, enzymes in medicines, therapeutic proteins, molecular and synthetic biology, expression and production of proteins of interest in the pharmaceutical industry."
$d.Paragraphs.Item(9).Range.Text = "Fornecer conhecimentos básicos aos estudantes do curso de Engenharia Bioquímica sobre os aspectos moleculares da utilização de microrganismos na obtenção novos compostos e na produção de moléculas de interesse farmacêutico (fármacos, insumos e de diagnóstico)."
$d.Paragraphs.Item(11).Range.Text = "1. Fundamentos de biotecnologia moderna;" + $nl + "2. Biotecnologia voltada a terapêutica;" + $nl + "3. Biologia molecular voltada a indústria farmacêutica;" + $nl + "4. Produção microbiana de agentes terapêuticos;" + $nl + "5. Tecnologias de alta eficiência para prospecção de novas moléculas;" + $nl + "6. Desenho racional de moléculas de interesse terapêutico;" + $nl + "7. Produção de agentes terapêuticos em microrganismos;"
$d.Paragraphs.Item(12).Range.Text = "Provide basic knowledge to Biochemical Engineering students on the molecular aspects of using microorganisms to obtain new compounds and produce molecules of pharmaceutical interest (drugs, inputs and diagnostics)."
$d.Paragraphs.Item(14).Range.Text = "A avaliação será feita por meio de provas escritas (P1 e P2)."
$d.Paragraphs.Item(19).Range.Text = "4873328 - Fernando Segato"

# --- Paragraph 17 (Avaliacao list): three content runs interleaved with bold
# labels (Metodo / Criterio / Norma de recuperacao) must be updated in place,
# leaving the labels untouched. Processed as a find/replace chain scoped to the
# paragraph, last-rotated-value first, so no search text is ambiguous.
$p17 = $d.Paragraphs.Item(17)

$r17c = $d.Range($p17.Range.Start, $p17.Range.End)
$r17c.Find.Execute("A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2", $true, $false, $false, $false, $false, $true, 1, $false, "1. Shayna Cox Gad (2007). Handbook of Pharmaceutical Biotechnology. John Wiley & Sons, New Jersey.^l2. Heinrich Klefenz (2002). Industrial Pharmaceutical Biotechnology. Wiley-VCH Verlag GmbH.^l3. Michael J. Groves (2006). Pharmaceutical Biotechnology. Taylor and Francis Group, USA.", 2) | Out-Null

$p17 = $d.Paragraphs.Item(17)
$r17b = $d.Range($p17.Range.Start, $p17.Range.End)
$r17b.Find.Execute("A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2", $true, $false, $false, $false, $false, $true, 1, $false, "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2", 2) | Out-Null

$p17 = $d.Paragraphs.Item(17)
$r17a = $d.Range($p17.Range.Start, $p17.Range.End)
$r17a.Find.Execute("A avaliação será feita por meio de provas escritas (P1 e P2).", $true, $false, $false, $false, $false, $true, 1, $false, "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2", 2) | Out-Null

Write-Output "edit complete"

